$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 1st team member's entry (Jakub Majer row, columns B:D):
# file changes from "clasa1.cs" to "diagram i funkcjonalnosci.docx" and its line count from 100 to 25
$ws.Range("C7").Value = "diagram i funkcjonalnosci.docx"
$ws.Range("D7").Value = 25

# Update the 4th team member's entry (Szymon Czapla row, columns K:M):
# file changes from "clasa3.cs" to "diagram i funkcjonalnosci.docx" (line count of 25 stays the same)
$ws.Range("L7").Value = "diagram i funkcjonalnosci.docx"

# Move / update the active selection to match the author's last saved cursor position
$ws.Range("I14").Select() | Out-Null
